$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1318
$ws.Range("I43").Value = 999
$ws.Range("J43").Value = 1397.75
$ws.Range("K43").Value = 999
$ws.Range("L43").Value = 1397.75
$ws.Range("M43").Value = -930
$ws.Range("N43").Value = -1535.75

$ws.Range("H51").Value = 2073.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2073.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2073.5
$ws.Range("N51").Value = -3041.5

$ws.Range("H74").Value = 3550
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 3825
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 3825
$ws.Range("M74").Value = -2064
$ws.Range("N74").Value = -5697

$ws.Range("H77").Value = 3550
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 3825
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 19125
$ws.Range("M77").Value = -10320
$ws.Range("N77").Value = -28485

$ws.Range("H116").Value = 11020.077
$ws.Range("I116").Value = 16664.285
$ws.Range("J116").Value = 4435.1665
$ws.Range("K116").Value = 16664.285
$ws.Range("L116").Value = 4435.1665
$ws.Range("M116").Value = -13222.285
$ws.Range("N116").Value = -11319.1665

$ws.Range("H125").Value = 1195.3572
$ws.Range("I125").Value = 1157.091
$ws.Range("J125").Value = 1335.6666
$ws.Range("K125").Value = 10413.819
$ws.Range("L125").Value = 12020.9994
$ws.Range("M125").Value = -7953.819
$ws.Range("N125").Value = -16940.9994

$ws.Range("H137").Value = 2031.9231
$ws.Range("I137").Value = 1772.6
$ws.Range("J137").Value = 2194
$ws.Range("K137").Value = 5317.799999999999
$ws.Range("L137").Value = 6582
$ws.Range("M137").Value = -2767.799999999999
$ws.Range("N137").Value = -11682

$ws.Range("H138").Value = 1760.8226
$ws.Range("I138").Value = 1268.6786
$ws.Range("J138").Value = 2166.1177
$ws.Range("K138").Value = 3806.0358
$ws.Range("L138").Value = 6498.353099999999
$ws.Range("M138").Value = 1333.9642
$ws.Range("N138").Value = -16778.3531

$ws.Range("H141").Value = 5274.8887
$ws.Range("I141").Value = 3914
$ws.Range("J141").Value = 6976
$ws.Range("K141").Value = 11742
$ws.Range("L141").Value = 20928
$ws.Range("M141").Value = -6562
$ws.Range("N141").Value = -31288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 795.5
$ws.Range("I5").Value = 91
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 91
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 21
$ws.Range("N5").Value = -1724

$ws.Range("H32").Value = 3636.6667
$ws.Range("I32").Value = 2188.5593
$ws.Range("J32").Value = 12180.5
$ws.Range("K32").Value = 2188.5593
$ws.Range("L32").Value = 12180.5
$ws.Range("M32").Value = -1901.5593
$ws.Range("N32").Value = -12754.5

$ws.Range("H61").Value = 3954.6333
$ws.Range("I61").Value = 3031.353
$ws.Range("J61").Value = 5162
$ws.Range("K61").Value = 3031.353
$ws.Range("L61").Value = 5162
$ws.Range("M61").Value = -2819.353
$ws.Range("N61").Value = -5586

$ws.Range("H122").Value = 1317.7142
$ws.Range("I122").Value = 1337.3334
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 4012.0002
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1562.0002
$ws.Range("N122").Value = -8500

$ws.Range("H132").Value = 1433.1569
$ws.Range("I132").Value = 970.6053000000001
$ws.Range("J132").Value = 2785.2307
$ws.Range("K132").Value = 2911.8159
$ws.Range("L132").Value = 8355.6921
$ws.Range("M132").Value = -381.8159000000001
$ws.Range("N132").Value = -13415.6921

$ws.Range("H136").Value = 3954.6333
$ws.Range("I136").Value = 3031.353
$ws.Range("J136").Value = 5162
$ws.Range("K136").Value = 9094.059000000001
$ws.Range("L136").Value = 15486
$ws.Range("M136").Value = -6544.059000000001
$ws.Range("N136").Value = -20586

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 795.5
$ws.Range("I4").Value = 91
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 91
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = 24
$ws.Range("N4").Value = -1730

$ws.Range("H134").Value = 5433.6216
$ws.Range("I134").Value = 6295.037
$ws.Range("J134").Value = 3107.8
$ws.Range("K134").Value = 18885.111
$ws.Range("L134").Value = 9323.400000000001
$ws.Range("M134").Value = -16350.111
$ws.Range("N134").Value = -14393.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3066.3333
$ws.Range("I16").Value = 1200
$ws.Range("J16").Value = 3999.5
$ws.Range("K16").Value = 1200
$ws.Range("L16").Value = 3999.5
$ws.Range("M16").Value = -913
$ws.Range("N16").Value = -4573.5

$ws.Range("H22").Value = 876.6
$ws.Range("I22").Value = 372.5
$ws.Range("J22").Value = 1059.909
$ws.Range("K22").Value = 372.5
$ws.Range("L22").Value = 1059.909
$ws.Range("M22").Value = -22.5
$ws.Range("N22").Value = -1759.909

$ws.Range("H31").Value = 2326.8462
$ws.Range("I31").Value = 2081.6667
$ws.Range("J31").Value = 2456.647
$ws.Range("K31").Value = 2081.6667
$ws.Range("L31").Value = 2456.647
$ws.Range("M31").Value = -1786.6667
$ws.Range("N31").Value = -3046.647

$ws.Range("H34").Value = 2326.8462
$ws.Range("I34").Value = 2081.6667
$ws.Range("J34").Value = 2456.647
$ws.Range("K34").Value = 2081.6667
$ws.Range("L34").Value = 2456.647
$ws.Range("M34").Value = -1879.6667
$ws.Range("N34").Value = -2860.647

$ws.Range("H113").Value = 3066.3333
$ws.Range("I113").Value = 1200
$ws.Range("J113").Value = 3999.5
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 3999.5
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -8339.5

$ws.Range("H132").Value = 2389.5806
$ws.Range("I132").Value = 1555.7142
$ws.Range("J132").Value = 4140.7
$ws.Range("K132").Value = 4667.142599999999
$ws.Range("L132").Value = 12422.1
$ws.Range("M132").Value = -2137.142599999999
$ws.Range("N132").Value = -17482.1

$ws.Range("H134").Value = 1132.5714
$ws.Range("I134").Value = 1131.8108
$ws.Range("J134").Value = 1138.2
$ws.Range("K134").Value = 3395.4324
$ws.Range("L134").Value = 3414.6
$ws.Range("M134").Value = -860.4323999999997
$ws.Range("N134").Value = -8484.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 190
$ws.Range("I2").Value = 216.66667
$ws.Range("J2").Value = 110
$ws.Range("K2").Value = 1300.00002
$ws.Range("L2").Value = 660
$ws.Range("M2").Value = -1187.00002
$ws.Range("N2").Value = -886

$ws.Range("H82").Value = 2000
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 6000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -5594

$ws.Range("H85").Value = 2000
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 6000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -4596

$ws.Range("H131").Value = 839.9299999999999
$ws.Range("I131").Value = 625.4
$ws.Range("J131").Value = 851.2210700000001
$ws.Range("K131").Value = 1876.2
$ws.Range("L131").Value = 2553.66321
$ws.Range("M131").Value = 3163.8
$ws.Range("N131").Value = -12633.66321

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2240.9
$ws.Range("I97").Value = 2201.4285
$ws.Range("J97").Value = 2333
$ws.Range("K97").Value = 2201.4285
$ws.Range("L97").Value = 2333
$ws.Range("M97").Value = -1705.4285
$ws.Range("N97").Value = -3325

$ws.Range("H102").Value = 1777.7667
$ws.Range("I102").Value = 1797.6428
$ws.Range("J102").Value = 1499.5
$ws.Range("K102").Value = 1797.6428
$ws.Range("L102").Value = 1499.5
$ws.Range("M102").Value = -175.6428000000001
$ws.Range("N102").Value = -4743.5

$ws.Range("H122").Value = 2544
$ws.Range("I122").Value = 1899
$ws.Range("J122").Value = 2728.2856
$ws.Range("K122").Value = 5697
$ws.Range("L122").Value = 8184.8568
$ws.Range("M122").Value = -3247
$ws.Range("N122").Value = -13084.8568

$ws.Range("H127").Value = 33866
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 33866
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 33866
$ws.Range("N127").Value = -43786

$ws.Range("H132").Value = 1042200.3
$ws.Range("I132").Value = 1540835.1
$ws.Range("J132").Value = 3377.8333
$ws.Range("K132").Value = 4622505.300000001
$ws.Range("L132").Value = 10133.4999
$ws.Range("M132").Value = -4619975.300000001
$ws.Range("N132").Value = -15193.4999

$ws.Range("H141").Value = 30000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 30000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4270
$ws.Range("I22").Value = 5116.6665
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 5116.6665
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -4821.6665
$ws.Range("N22").Value = -3590

$ws.Range("H27").Value = 4270
$ws.Range("I27").Value = 5116.6665
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 5116.6665
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -5009.6665
$ws.Range("N27").Value = -3214

$ws.Range("H40").Value = 2912.6
$ws.Range("I40").Value = 2254.0908
$ws.Range("J40").Value = 4723.5
$ws.Range("K40").Value = 2254.0908
$ws.Range("L40").Value = 4723.5
$ws.Range("M40").Value = -2118.0908
$ws.Range("N40").Value = -4995.5

$ws.Range("H46").Value = 2600.25
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 3700.5
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 3700.5
$ws.Range("M46").Value = -1312
$ws.Range("N46").Value = -4076.5

$ws.Range("H55").Value = 292.48572
$ws.Range("I55").Value = 250.14815
$ws.Range("J55").Value = 435.375
$ws.Range("K55").Value = 250.14815
$ws.Range("L55").Value = 435.375
$ws.Range("M55").Value = -77.14814999999999
$ws.Range("N55").Value = -781.375

$ws.Range("H132").Value = 3036.8215
$ws.Range("I132").Value = 2435.2222
$ws.Range("J132").Value = 4119.7
$ws.Range("K132").Value = 7305.6666
$ws.Range("L132").Value = 12359.1
$ws.Range("M132").Value = -4775.6666
$ws.Range("N132").Value = -17419.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 56687.355
$ws.Range("I122").Value = 56687.355
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 170062.065
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -167612.065

$ws.Range("H132").Value = 1654.0264
$ws.Range("I132").Value = 1072.8334
$ws.Range("J132").Value = 3833.5
$ws.Range("K132").Value = 3218.5002
$ws.Range("L132").Value = 11500.5
$ws.Range("M132").Value = -688.5001999999999
$ws.Range("N132").Value = -16560.5
